# Applies the "Screen Flash and finish chamer and flag added" edit:
#  1) Spell-check style run-splitting (w:proofErr spellStart/spellEnd) around
#     the made-up / foreign words in the intro paragraphs (text unchanged).
#  2) Re-worked "To Do:" list: several items rewritten, the paragraph-mark
#     language formatting (w:lang en-US) stripped from the surviving items,
#     three now-redundant bullets removed, and the last bullet gets the same
#     spell-check run-splitting treatment.
#
# Implemented with Range.InsertXML, which replaces the full contents of the
# target Range (here: a whole paragraph, mark included) with literal OOXML,
# and Range.Delete(), which removes a whole paragraph (mark included) and
# merges its neighbours. Operations are issued in strictly descending
# paragraph-index order so that earlier (lower) indices stay valid while
# later ones shift from the deletions.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 19: "-Nochmals Trackobjekte ... erweiterbarkeit)" ----------
$xml19 = '<w:p ' + $wns + '>' +
  '<w:r><w:t xml:space="preserve">-Nochmals </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Trackobjekte</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> automatisch einscannen probieren</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> (dann im </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>ggd</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> erwähnen bezüglich einfacher </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>erweiterbarkeit</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>)</w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(19).Range.InsertXML($xml19)

# --- Paragraph 18: "-Sound" (drop the en-US paragraph-mark formatting) ----
$xml18 = '<w:p ' + $wns + '><w:r><w:t>-Sound</w:t></w:r></w:p>'
$d.Paragraphs(18).Range.InsertXML($xml18)

# --- Paragraph 17: "-Fähnchen für das Ende" -> removed entirely -----------
$d.Paragraphs(17).Range.Delete()

# --- Paragraph 16: "-Explosion Element einbauen." (old slot) -> removed ---
$d.Paragraphs(16).Range.Delete()

# --- Paragraph 15: "-Fix End Screen (...)" (old slot) -> removed ----------
$d.Paragraphs(15).Range.Delete()

# --- Paragraph 14: "-Fall of Map Message oder Effekt" + bookmark ----------
#     becomes "-Explosion Element einbauen." (bookmark kept, lang dropped)
$xml14 = '<w:p ' + $wns + '>' +
  '<w:r><w:t>-Explosion Element einbauen.</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '</w:p>'
$d.Paragraphs(14).Range.InsertXML($xml14)

# --- Paragraph 13: "-Welcome Screen Mausrad" -------------------------------
#     becomes "-Fix End Screen (Text, Replay Button)" (keeps en-US lang)
$xml13 = '<w:p ' + $wns + '>' +
  '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>-Fix End Screen (Text, Replay Button)</w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(13).Range.InsertXML($xml13)

# --- Paragraph 5: "-blah einfache trackerweiterung" (proofErr splitting) --
$xml5 = '<w:p ' + $wns + '>' +
  '<w:r><w:t>-</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>blah</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> einfache </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>trackerweiterung</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
$d.Paragraphs(5).Range.InsertXML($xml5)

# --- Paragraph 4: "-blah Volumentrigger" (proofErr splitting) -------------
$xml4 = '<w:p ' + $wns + '>' +
  '<w:r><w:t>-</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>blah</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> Volumentrigger</w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(4).Range.InsertXML($xml4)

# --- Paragraph 3: "Kurze Beschreibung des Spiels ..." (proofErr splitting) --
$xml3 = '<w:p ' + $wns + '>' +
  '<w:r><w:t>Kurze Beschreibung des Spiels</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>blah</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>spielobjekte</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>blah</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> warum ich normale </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>gravitation</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> und </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>rollmode</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> gemacht habe</w:t></w:r>' +
  '</w:p>'
$d.Paragraphs(3).Range.InsertXML($xml3)
